$d = $word.ActiveDocument

function FindParaByText($needle) {
    for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
        $pp = $d.Paragraphs.Item($i)
        if ($pp.Range.Text -like $needle) {
            return $pp
        }
    }
    return $null
}

# Paste `sampleRange` (a single-character donor range whose rPr is the
# formatting we want to clone) at $cursor, then overwrite the pasted
# character(s) with $text -- this keeps the donor's run properties
# (sz/szCs/u/lang/...) intact while giving us the text we actually want.
# Returns the cursor position right after the inserted text.
function PasteRun($cursor, $sampleRange, $text) {
    $len = $sampleRange.Text.Length
    $sampleRange.Copy()
    $ip = $d.Range($cursor, $cursor)
    $ip.Paste()
    $runRange = $d.Range($cursor, $cursor + $len)
    $runRange.Text = $text
    return $cursor + $text.Length
}

# ------------------------------------------------------------------
# Step 1: heading text "ATRIBUICOES" -> "RESPONSABILIDADE POR FUNCAO"
# ------------------------------------------------------------------
$d.Content.Find.Execute("ATRIBUIÇÕES", $true, $false, $false, $false, $false, $true, 1, $false, "RESPONSABILIDADE POR FUNÇÃO", 2) | Out-Null

# ------------------------------------------------------------------
# grab formatting donors *before* we start rewriting paragraphs
# ------------------------------------------------------------------
# sz=24 (no szCs) + lang=pt-PT  -> lives only on the "Cabe a Chefe..." runs
$cabe = FindParaByText("*Cabe a Chefe*")
$noCsSample = $d.Range($cabe.Range.Start + 3, $cabe.Range.Start + 4)
Write-Output ("noCsSample=[" + $noCsSample.Text + "]")

# sz=24 + szCs=24 + lang=pt-PT
$szCsSrcPara = FindParaByText("*CELOG – Centro Logístico*")
$szCsSample = $d.Range($szCsSrcPara.Range.Start, $szCsSrcPara.Range.Start + 1)
Write-Output ("szCsSample=[" + $szCsSample.Text + "]")

# sz=24 + szCs=24 + u=single + lang=pt-PT  (section-heading style)
$headingSrcPara = FindParaByText("RESPONSABILIDADE POR FUNÇÃO*")
$headingSample = $d.Range($headingSrcPara.Range.Start, $headingSrcPara.Range.Start + 1)
Write-Output ("headingSample=[" + $headingSample.Text + "]")

# ------------------------------------------------------------------
# Step 2: rewrite the "Cabe a Chefe..." paragraph
# ------------------------------------------------------------------
$cabe = FindParaByText("*Cabe a Chefe*")
$pr = $cabe.Range
$pr.End = $pr.End - 1
$cursor = $pr.Start
$pr.Text = ""

$cursor = PasteRun $cursor $szCsSample "As "
$cursor = PasteRun $cursor $szCsSample "responsabilidades por função são delineadas sinteticamente no Regimento Interno do CELOG e as atividades relacionados aos seus cumprimentos estão contempladas n"
$cursor = PasteRun $cursor $noCsSample "os seguintes processos:"

$cabe = FindParaByText("*seguintes processos:*")
Write-Output ("Step2 paragraph now: [" + $cabe.Range.Text + "]")

# ------------------------------------------------------------------
# Step 3: insert a blank paragraph + two new headed paragraphs after
# "PLOG0009 - Elaboracao de plano de inspecao"
# ------------------------------------------------------------------
$plog9 = FindParaByText("*plano de inspeção*")
$plog9Range = $plog9.Range
$plog9Len = $plog9Range.Text.Length
$plog9Range.Copy()

$pasteAt = $plog9.Range.End
$ip = $d.Range($pasteAt, $pasteAt)
$ip.Paste()
$ip2 = $d.Range($pasteAt + $plog9Len, $pasteAt + $plog9Len)
$ip2.Paste()
$ip3 = $d.Range($pasteAt + 2 * $plog9Len, $pasteAt + 2 * $plog9Len)
$ip3.Paste()

# Re-fetch paragraph objects for the three freshly pasted paragraphs
$plog9 = FindParaByText("*plano de inspeção*")
$blankPara = $plog9.Next()
$headingPara = $blankPara.Next()
$bodyPara = $headingPara.Next()

Write-Output ("blank=[" + $blankPara.Range.Text + "]")
Write-Output ("heading=[" + $headingPara.Range.Text + "]")
Write-Output ("body=[" + $bodyPara.Range.Text + "]")

# --- blank paragraph: strip numbering & text ---
$blankPara.Range.ListFormat.RemoveNumbers()
$blankPara = $plog9.Next()
$br = $blankPara.Range
$br.End = $br.End - 1
$br.Text = ""

# --- heading paragraph: ilvl 1, clear text, insert underlined run ---
$headingPara = $blankPara.Next()
$headingPara.Range.ListFormat.ListLevelNumber = 2
$headingPara = $blankPara.Next()
$hr = $headingPara.Range
$hr.End = $hr.End - 1
$hcursor = $hr.Start
$hr.Text = ""
$hcursor = PasteRun $hcursor $headingSample "RESPONSABILIDADE NO INTER-RELACIONAMENTO ENTRE OS SETORES"

# --- body paragraph: ilvl 2 (unchanged, but re-assert it), clear text, insert two runs ---
$headingPara = $blankPara.Next()
$bodyPara = $headingPara.Next()
$bodyPara.Range.ListFormat.ListLevelNumber = 3
$bodyPara = $headingPara.Next()
$bro = $bodyPara.Range
$bro.End = $bro.End - 1
$bcursor = $bro.Start
$bro.Text = ""

$bcursor = PasteRun $bcursor $szCsSample "As responsa"
$bcursor = PasteRun $bcursor $szCsSample "bilidades no inter-relacionamento entre setores são apresentadas detalhadamente nos PLOG relacionados no item 2.2 desta NPA."

$blankPara = $plog9.Next()
$headingPara = $blankPara.Next()
$bodyPara = $headingPara.Next()
Write-Output ("Final blank=[" + $blankPara.Range.Text + "]")
Write-Output ("Final heading=[" + $headingPara.Range.Text + "]")
Write-Output ("Final body=[" + $bodyPara.Range.Text + "]")
Write-Output "DONE"
